# "Se procesan de nuevo los datos con las nuevas dimensiones curadas"
#
# Five fields (cnae-ano, seccion-1-letra-descripcion, mes-nombre,
# clasificacion, sexo) -- columns E, G, I, L, M -- are recurated from
# "dimension" to "measure":
#   row 2: the iaest-dimension:<field> URN becomes iaest-measure:<field>
#   row 3: the "dim" marker becomes "medida"
#   row 4: the "skos:Concept" type becomes "xsd:int"
#   row 5: the per-field mapping-<field>.xlsx reference is dropped entirely
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$recuratedCols = "E", "G", "I", "L", "M"

$row2New = @{
    "E" = "iaest-measure:cnae-ano"
    "G" = "iaest-measure:seccion-1-letra-descripcion"
    "I" = "iaest-measure:mes-nombre"
    "L" = "iaest-measure:clasificacion"
    "M" = "iaest-measure:sexo"
}

foreach ($col in $recuratedCols) {
    $ws.Range($col + "2").Value = $row2New[$col]
    $ws.Range($col + "3").Value = "medida"
    $ws.Range($col + "4").Value = "xsd:int"
    # The mapping-file row no longer applies to measures, so the cell is
    # removed entirely (not just blanked).
    $ws.Range($col + "5").Clear()
}
